{"js": "// Update the worksheet date and the 25 three-digit-by-one-digit division\n// problems in the table, keeping every run's original formatting intact.\n//\n// The body is: async (context) => { ... } (this is the body text).\n\nconst body = context.document.body;\n\n// ---- 1. Header date paragraph --------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text.trim() === \"2025-08-30 Saturday\") {\n  dateParagraph.getRange(\"Content\").insertText(\"2025-08-31 Sunday\", \"Replace\");\n  await context.sync();\n}\n\n// ---- 2. Division problems inside the table --------------------------------\n// The problems live in table rows 0, 4, 8, 12, 16 (5 columns each); the rows\n// in between hold the students' blank answer space. Addressing cells by\n// (row, column) avoids any ambiguity from the repeated \"855\u00f79=\" problem.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  { row: 0, col: 0, oldText: \"269\u00f75=\", newText: \"477\u00f75=\" },\n  { row: 0, col: 1, oldText: \"748\u00f74=\", newText: \"869\u00f76=\" },\n  { row: 0, col: 2, oldText: \"823\u00f74=\", newText: \"837\u00f72=\" },\n  { row: 0, col: 3, oldText: \"985\u00f77=\", newText: \"935\u00f72=\" },\n  { row: 0, col: 4, oldText: \"882\u00f75=\", newText: \"535\u00f74=\" },\n\n  { row: 4, col: 0, oldText: \"624\u00f72=\", newText: \"702\u00f78=\" },\n  { row: 4, col: 1, oldText: \"595\u00f72=\", newText: \"505\u00f75=\" },\n  { row: 4, col: 2, oldText: \"779\u00f73=\", newText: \"777\u00f76=\" },\n  { row: 4, col: 3, oldText: \"889\u00f77=\", newText: \"420\u00f79=\" },\n  { row: 4, col: 4, oldText: \"979\u00f73=\", newText: \"183\u00f74=\" },\n\n  { row: 8, col: 0, oldText: \"855\u00f79=\", newText: \"661\u00f76=\" },\n  { row: 8, col: 1, oldText: \"104\u00f77=\", newText: \"453\u00f76=\" },\n  { row: 8, col: 2, oldText: \"388\u00f72=\", newText: \"635\u00f79=\" },\n  { row: 8, col: 3, oldText: \"391\u00f72=\", newText: \"672\u00f77=\" },\n  { row: 8, col: 4, oldText: \"395\u00f79=\", newText: \"489\u00f78=\" },\n\n  { row: 12, col: 0, oldText: \"946\u00f72=\", newText: \"730\u00f75=\" },\n  { row: 12, col: 1, oldText: \"825\u00f72=\", newText: \"916\u00f79=\" },\n  { row: 12, col: 2, oldText: \"690\u00f79=\", newText: \"868\u00f78=\" },\n  { row: 12, col: 3, oldText: \"995\u00f77=\", newText: \"287\u00f76=\" },\n  { row: 12, col: 4, oldText: \"855\u00f79=\", newText: \"152\u00f73=\" },\n\n  { row: 16, col: 0, oldText: \"503\u00f76=\", newText: \"693\u00f79=\" },\n  { row: 16, col: 1, oldText: \"202\u00f72=\", newText: \"514\u00f78=\" },\n  { row: 16, col: 2, oldText: \"546\u00f78=\", newText: \"373\u00f77=\" },\n  { row: 16, col: 3, oldText: \"329\u00f79=\", newText: \"519\u00f72=\" },\n  { row: 16, col: 4, oldText: \"591\u00f78=\", newText: \"811\u00f74=\" }\n];\n\nconst cells = updates.map((u) => table.getCell(u.row, u.col));\ncells.forEach((cell) => cell.body.load(\"text\"));\nawait context.sync();\n\nfor (let i = 0; i < updates.length; i++) {\n  const u = updates[i];\n  const cell = cells[i];\n  if (cell.body.text.trim() === u.oldText) {\n    cell.body.getRange(\"Content\").insertText(u.newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 three-digit-by-one-digit division\n# problems in the table, keeping every run's original formatting intact.\n#\n# $d.Cell(row, col).Range.Text = \"...\" rewrites only the cell's text content\n# (the COM layer keeps the existing run/paragraph formatting and the\n# trailing end-of-cell marker), exactly like typing over the selected text\n# in Word.\n\n$d = $word.ActiveDocument\n\n# ---- 1. Header date paragraph ---------------------------------------------\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text.TrimEnd(\"`r\") -eq \"2025-08-30 Saturday\") {\n    $dateParagraph.Range.Text = \"2025-08-31 Sunday\"\n}\n\n# ---- 2. Division problems inside the table ---------------------------------\n# The problems live in table rows 1, 5, 9, 13, 17 (1-based; 5 columns each);\n# the rows in between hold the students' blank answer space. Addressing\n# cells by (row, column) avoids any ambiguity from the repeated \"855\u00f79=\"\n# problem.\n$table = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 1;  Col = 1; Old = \"269\u00f75=\"; New = \"477\u00f75=\" },\n    @{ Row = 1;  Col = 2; Old = \"748\u00f74=\"; New = \"869\u00f76=\" },\n    @{ Row = 1;  Col = 3; Old = \"823\u00f74=\"; New = \"837\u00f72=\" },\n    @{ Row = 1;  Col = 4; Old = \"985\u00f77=\"; New = \"935\u00f72=\" },\n    @{ Row = 1;  Col = 5; Old = \"882\u00f75=\"; New = \"535\u00f74=\" },\n\n    @{ Row = 5;  Col = 1; Old = \"624\u00f72=\"; New = \"702\u00f78=\" },\n    @{ Row = 5;  Col = 2; Old = \"595\u00f72=\"; New = \"505\u00f75=\" },\n    @{ Row = 5;  Col = 3; Old = \"779\u00f73=\"; New = \"777\u00f76=\" },\n    @{ Row = 5;  Col = 4; Old = \"889\u00f77=\"; New = \"420\u00f79=\" },\n    @{ Row = 5;  Col = 5; Old = \"979\u00f73=\"; New = \"183\u00f74=\" },\n\n    @{ Row = 9;  Col = 1; Old = \"855\u00f79=\"; New = \"661\u00f76=\" },\n    @{ Row = 9;  Col = 2; Old = \"104\u00f77=\"; New = \"453\u00f76=\" },\n    @{ Row = 9;  Col = 3; Old = \"388\u00f72=\"; New = \"635\u00f79=\" },\n    @{ Row = 9;  Col = 4; Old = \"391\u00f72=\"; New = \"672\u00f77=\" },\n    @{ Row = 9;  Col = 5; Old = \"395\u00f79=\"; New = \"489\u00f78=\" },\n\n    @{ Row = 13; Col = 1; Old = \"946\u00f72=\"; New = \"730\u00f75=\" },\n    @{ Row = 13; Col = 2; Old = \"825\u00f72=\"; New = \"916\u00f79=\" },\n    @{ Row = 13; Col = 3; Old = \"690\u00f79=\"; New = \"868\u00f78=\" },\n    @{ Row = 13; Col = 4; Old = \"995\u00f77=\"; New = \"287\u00f76=\" },\n    @{ Row = 13; Col = 5; Old = \"855\u00f79=\"; New = \"152\u00f73=\" },\n\n    @{ Row = 17; Col = 1; Old = \"503\u00f76=\"; New = \"693\u00f79=\" },\n    @{ Row = 17; Col = 2; Old = \"202\u00f72=\"; New = \"514\u00f78=\" },\n    @{ Row = 17; Col = 3; Old = \"546\u00f78=\"; New = \"373\u00f77=\" },\n    @{ Row = 17; Col = 4; Old = \"329\u00f79=\"; New = \"519\u00f72=\" },\n    @{ Row = 17; Col = 5; Old = \"591\u00f78=\"; New = \"811\u00f74=\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $table.Cell($u.Row, $u.Col)\n    $current = $cell.Range.Text.TrimEnd(\"`r\", [char]7)\n    if ($current -eq $u.Old) {\n        $cell.Range.Text = $u.New\n    }\n}\n"}
